# mainconfig.xlsx - "fixed tests missing new keyword in config files"
# Add a new "MPI executable prefix" row to the MAIN Config. sheet, between
# the existing "MPI tasks" row and the "Batch system" row.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "MAIN Config." ---------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row 13 (pushes old row 13 -> 14, old row 14 -> 15) and give
# it the same formatting as the row above it (row 12, "MPI tasks").
$ws1.Rows.Item(13).Insert() | Out-Null
$ws1.Range("A12:B12").Copy() | Out-Null
$ws1.Range("A13:B13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Cells.Item(13, 1).Value = "MPI executable prefix"
$ws1.Cells.Item(13, 2).Value = ""

# --- Active sheet / selection ---------------------------------------------
# The saved workbook now re-opens with "MAIN Config." as the active sheet
# (instead of "Experimental benchmarks"), with B33 selected there.
$ws1.Activate() | Out-Null
$ws1.Range("B33").Select() | Out-Null
